$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates: force Text format per-cell so the numeric-looking
# strings are stored as text (matching the original inlineStr cell type) ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.839.81'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.024.67'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.13'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.10'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.021.14'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.134'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000224'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.77'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.120'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.524.26'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.017.42'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.843.80'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.43'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '434.73'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.63'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.58'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.19'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.48'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.63'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.01'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0670'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.48'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0361'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '388.77'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.111'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.672.03'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.237'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.34'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.77'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.135'

# --- Column E (Volume(1h)) updates (already non-numeric text, e.g. "  -4.68%  ") ---
$ws.Range("E2").Value = '  -4.68%  '
$ws.Range("E3").Value = '  -5.94%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("E5").Value = '  -2.27%  '
$ws.Range("E6").Value = '  -7.67%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -5.97%  '
$ws.Range("E9").Value = '  -2.90%  '
$ws.Range("E10").Value = '  -6.92%  '
$ws.Range("E11").Value = '  -2.78%  '
$ws.Range("E12").Value = '  -3.46%  '
$ws.Range("E13").Value = '  -6.95%  '
$ws.Range("E14").Value = '  -7.48%  '
$ws.Range("E15").Value = '  +0.43%  '
$ws.Range("E16").Value = '  -5.94%  '
$ws.Range("E17").Value = '  -6.25%  '
$ws.Range("E18").Value = '  -4.72%  '
$ws.Range("E19").Value = '  -2.48%  '
$ws.Range("E20").Value = '  -6.60%  '
$ws.Range("E21").Value = '  -6.58%  '
$ws.Range("E22").Value = '  -5.19%  '
$ws.Range("E23").Value = '  -8.58%  '
$ws.Range("E24").Value = '  -4.40%  '
$ws.Range("E25").Value = '  -4.70%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("E28").Value = '  -4.42%  '
$ws.Range("E29").Value = '  -6.70%  '
$ws.Range("E30").Value = '  -7.96%  '
$ws.Range("E31").Value = '  -9.88%  '
$ws.Range("E32").Value = '  -7.64%  '
$ws.Range("E33").Value = '  -9.20%  '
$ws.Range("E34").Value = '  -12.22%  '
$ws.Range("E35").Value = '  -7.78%  '
$ws.Range("E36").Value = '  -4.90%  '
$ws.Range("E37").Value = '  -1.34%  '
$ws.Range("E38").Value = '  -8.87%  '
$ws.Range("E39").Value = '  +4.04%  '
$ws.Range("E40").Value = '  -8.52%  '
$ws.Range("E41").Value = '  -3.93%  '
$ws.Range("E42").Value = '  -2.53%  '
$ws.Range("E43").Value = '  -9.38%  '
$ws.Range("E44").Value = '  -5.89%  '
$ws.Range("E46").Value = '  -7.70%  '
$ws.Range("E47").Value = '  -6.37%  '
$ws.Range("E48").Value = '  -7.67%  '
$ws.Range("E49").Value = '  -3.81%  '
$ws.Range("E50").Value = '  -8.02%  '
$ws.Range("E51").Value = '  +2.17%  '
